$p = $ppt.ActivePresentation

# ------------------------------------------------------------------
# 1. Unhide the hidden slides (slides 1-5 had show="0" / Hidden=True).
#    Slide 6 was already visible and stays untouched.
# ------------------------------------------------------------------
for ($i = 1; $i -le 5; $i++) {
    $slide = $p.Slides.Item($i)
    $slide.SlideShowTransition.Hidden = 0
}

# ------------------------------------------------------------------
# 2. Fix the click-through (build) order of the four derivation
#    textboxes on slide 3 so they appear left-to-right instead of
#    right-to-left: spid 55/46/44/39 -> 39/44/46/55.
# ------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$seq = $s3.TimeLine.MainSequence
$shapes3 = $s3.Shapes

function GetShapeById($shapes, $id) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        if ($shapes.Item($i).Id -eq $id) { return $shapes.Item($i) }
    }
    return $null
}

# Capture the existing click-effect order (10 entries) before touching it.
$originalOrder = @()
for ($i = 1; $i -le $seq.Count; $i++) {
    $originalOrder += , $seq.Item($i).Shape.Id
}

# Desired order: the first four targets reversed, the rest unchanged.
$swap = @{ 55 = 39; 46 = 44; 44 = 46; 39 = 55 }
$newOrder = @()
foreach ($id in $originalOrder) {
    if ($swap.ContainsKey($id)) {
        $newOrder += , $swap[$id]
    } else {
        $newOrder += , $id
    }
}

# Recreate every click effect in the corrected order (PowerPoint's
# object model has no supported way to re-point an existing Effect at
# a different shape, so rebuild the sequence from scratch).
while ($seq.Count -gt 0) {
    $seq.Item(1).Delete()
}
foreach ($id in $newOrder) {
    $sh = GetShapeById $shapes3 $id
    $seq.AddEffect($sh, 1) | Out-Null
}
